$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header for column C
$ws.Range("C1").Value = "Full name"

# New row 2 data (existing rows 1-2 keep username/Password/danielUde/Password123@)
$ws.Range("C2").Value = "Ude Daniel"

# New row 3: Vandal Savage's credentials
$ws.Range("A3").Value = "Vsavage"
$ws.Range("B3").Value = "Savage123@"
$ws.Range("C3").Value = "Vandal Savage"

# Set column C width similar to the other data columns
$ws.Columns.Item(3).ColumnWidth = 13.85546875

# Move selection to D3 (next empty cell after the new data), matching Excel's
# natural cursor position after typing the last entry
$ws.Range("D3").Select()
